$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 59

# Columns A and D contain text that looks like a date / number
# ("2023-06-19" and "25"). Force them to be stored as text (matching
# the rest of the sheet, which uses inline strings for these columns)
# instead of letting Excel auto-convert them to a date serial / number.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-19"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "18:34:37"
$ws.Cells.Item($row, 3).Value = "Monday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "25"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 121997
$ws.Cells.Item($row, 6).Value = 133819
$ws.Cells.Item($row, 7).Value = 162346
$ws.Cells.Item($row, 8).Value = 133132
$ws.Cells.Item($row, 9).Value = 177181
$ws.Cells.Item($row, 10).Value = 114647
$ws.Cells.Item($row, 11).Value = 201533
$ws.Cells.Item($row, 12).Value = 225192
$ws.Cells.Item($row, 13).Value = 175482
$ws.Cells.Item($row, 14).Value = 103790
$ws.Cells.Item($row, 15).Value = 39181
$ws.Cells.Item($row, 16).Value = 33941
$ws.Cells.Item($row, 17).Value = 51809
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36434
$ws.Cells.Item($row, 20).Value = -1
